$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Column D holds numeric-looking price strings that are stored as
# plain text in the workbook (e.g. "67.572.90", "613.16"). A leading
# apostrophe forces Excel to keep the assigned value as literal text
# instead of silently re-parsing/rounding it as a number.

$ws.Range("D2").Value = "'67.572.90"
$ws.Range("E2").Value = "  -2.74%  "

$ws.Range("D3").Value = "'3.536.60"

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'613.16"
$ws.Range("E5").Value = "  -5.15%  "

$ws.Range("D6").Value = "'154.18"
$ws.Range("E6").Value = "  -3.39%  "

$ws.Range("D7").Value = "'3.531.95"
$ws.Range("E7").Value = "  -3.69%  "

$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("E9").Value = "  -2.06%  "

$ws.Range("E10").Value = "  -2.59%  "

$ws.Range("D11").Value = "'6.90"
$ws.Range("E11").Value = "  -2.78%  "

$ws.Range("E12").Value = "  -1.78%  "

$ws.Range("E13").Value = "  -3.47%  "

$ws.Range("D14").Value = "'32.21"
$ws.Range("E14").Value = "  -1.26%  "

$ws.Range("D15").Value = "'4.133.10"
$ws.Range("E15").Value = "  -3.68%  "

$ws.Range("D16").Value = "'3.541.97"
$ws.Range("E16").Value = "  -3.49%  "

$ws.Range("D17").Value = "'67.555.72"

$ws.Range("E18").Value = "  +0.57%  "

$ws.Range("E19").Value = "  -0.92%  "

$ws.Range("D20").Value = "'15.56"
$ws.Range("E20").Value = "  -2.60%  "

$ws.Range("D21").Value = "'454.45"
$ws.Range("E21").Value = "  -2.22%  "

$ws.Range("D22").Value = "'9.39"
$ws.Range("E22").Value = "  -3.84%  "

$ws.Range("D23").Value = "'0.643"
$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("D24").Value = "'78.62"
$ws.Range("E24").Value = "  -1.07%  "

$ws.Range("D25").Value = "'3.676.68"
$ws.Range("E25").Value = "  -3.72%  "

$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("D27").Value = "'0.0000122"
$ws.Range("E27").Value = "  -2.24%  "

$ws.Range("D28").Value = "'10.48"
$ws.Range("E28").Value = "  -2.64%  "

$ws.Range("E29").Value = "  -6.02%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'2.57"
$ws.Range("E30").Value = "  -1.61%  "

$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "'1.70"
$ws.Range("E31").Value = "  +1.77%  "

$ws.Range("D32").Value = "'1.01"
$ws.Range("E32").Value = "  +0.60%  "

$ws.Range("E33").Value = "  -3.84%  "

$ws.Range("E34").Value = "  -2.12%  "

$ws.Range("E35").Value = "  -3.53%  "

$ws.Range("E36").Value = "  -4.04%  "

$ws.Range("D37").Value = "'3.532.37"
$ws.Range("E37").Value = "  -3.53%  "

$ws.Range("D38").Value = "'8.01"
$ws.Range("E38").Value = "  -4.17%  "

$ws.Range("E39").Value = "  -0.03%  "

$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  -0.06%  "

$ws.Range("D41").Value = "'172.44"
$ws.Range("E41").Value = "  -3.30%  "

$ws.Range("E42").Value = "  -4.62%  "

$ws.Range("E43").Value = "  -1.49%  "

$ws.Range("D44").Value = "'2.11"
$ws.Range("E44").Value = "  -3.13%  "

$ws.Range("D45").Value = "'0.894"
$ws.Range("E45").Value = "  -3.51%  "

$ws.Range("D46").Value = "'29.51"
$ws.Range("E46").Value = "  +9.80%  "

$ws.Range("D47").Value = "'45.76"
$ws.Range("E47").Value = "  -1.86%  "

$ws.Range("E48").Value = "  -2.22%  "

$ws.Range("E49").Value = "  -2.09%  "

$ws.Range("E50").Value = "  -1.77%  "

$ws.Range("D51").Value = "'1.03"
$ws.Range("E51").Value = "  -2.11%  "

